$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.289852738380432
$ws.Range("B1").Value = 3.258122444152832
$ws.Range("C1").Value = 5.820755958557129
$ws.Range("D1").Value = 1.750788927078247
$ws.Range("E1").Value = 1.026689171791077
